$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (13-18) continuing the Test Plan table with
# Role/Responsibilities, Schedule, Risks, Entry/Exit Criteria, Approval.

$rows = @(
    @{ Num = 12; Section = "Role and Responsibilities"; Details = "Govind: UI; Partner: Backend" },
    @{ Num = 13; Section = "`tSchedule"; Details = "Writing: 7 Apr, Execution: 8–9 Apr, Bug Fixing: 10 Apr, Final: 11 Apr" },
    @{ Num = 14; Section = "Risks and Mitigation"; Details = "Internet issue → offline tools, Backend delay → mock data" },
    @{ Num = 15; Section = "`tEntry Criteria"; Details = "Local deployment, Working modules" },
    @{ Num = 16; Section = "Exit Criteria"; Details = "Major bugs fixed, All test cases passed" },
    @{ Num = 17; Section = "Approval"; Details = "`t(Leave blank or write “Faculty signature here”)" }
)

$rowIndex = 13
foreach ($r in $rows) {
    $ws.Range("A$rowIndex").Value = $r.Num
    $ws.Range("B$rowIndex").Value = $r.Section
    $ws.Range("C$rowIndex").Value = $r.Details

    $rng = $ws.Range("A" + $rowIndex + ":C" + $rowIndex)
    $rng.HorizontalAlignment = -4108  # xlCenter
    $rng.VerticalAlignment = -4160    # xlTop
    $rng.WrapText = $true

    $rowIndex++
}

# Match the wrapped-text row heights Excel computed for the new rows.
$ws.Rows(13).RowHeight = 28.8
$ws.Rows(14).RowHeight = 57.6
$ws.Rows(15).RowHeight = 43.2
$ws.Rows(16).RowHeight = 28.8
$ws.Rows(17).RowHeight = 28.8
$ws.Rows(18).RowHeight = 43.2

$ws.Range("C19").Select()
